# Logbook update: append the client confirmation + screenshot filename
# reference to the "Screenshot multichat bersama kelompok client" bullet
# in the first Logbook table.

$d = $word.ActiveDocument

$r = $d.Content
$found = $r.Find.Execute(
    "Screenshot multichat bersama kelompok client",
    $true,   # MatchCase
    $true,   # MatchWholeWord
    $false,  # MatchWildcards
    $false,  # MatchSoundsLike
    $false,  # MatchAllWordForms
    $true,   # Forward
    1,       # Wrap = wdFindContinue
    $false,  # Format
    "",      # ReplaceWith
    0        # Replace = wdReplaceNone
)

if ($found) {
    # Move to the very end of the matched text, right after "...client".
    $r.Collapse(0)

    # Three separate inserts, mirroring the three new runs added by the
    # original edit: " : ", the opening curly quote, and the closing
    # quoted filename. Each inherits the run formatting (Times New Roman /
    # en-US) already in effect at the insertion point.
    $r.InsertAfter(" : ")
    $r.Collapse(0)

    $r.InsertAfter([char]8220)
    $r.Collapse(0)

    $r.InsertAfter("SS multichat.jpg" + [char]8221)
    $r.Collapse(0)
}
